$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.950.44'
$ws.Range("E2").Value = '  -1.46%  '
$ws.Range("D3").Value = '1.639.18'
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("E8").Value = '  -0.79%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.64'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0796'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.25%  '
$ws.Range("D12").Value = '1.865.82'
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.28'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = '1.636.61'
$ws.Range("E14").Value = '  -1.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.544'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.01'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.98%  '
$ws.Range("D18").Value = '25.953.67'
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.10%  '
$ws.Range("E21").Value = '  -1.85%  '
$ws.Range("E22").Value = '  -1.59%  '
$ws.Range("E23").Value = '  -0.88%  '
$ws.Range("E24").Value = '  +0.41%  '
$ws.Range("B25").Value = 'Stellar'
$ws.Range("C25").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.130'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.86%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '143.40'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("B27").Value = 'BinanceUSD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.32%  '
$ws.Range("E28").Value = '  -1.90%  '
$ws.Range("E29").Value = '  -0.53%  '
$ws.Range("E30").Value = '  -0.71%  '
$ws.Range("E31").Value = '  -1.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.24'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.53'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.11%  '
$ws.Range("E35").Value = '  +1.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.902'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.69%  '
$ws.Range("D37").Value = '1.135.68'
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("E38").Value = '  -2.19%  '
$ws.Range("E39").Value = '  -1.43%  '
$ws.Range("E40").Value = '  -0.29%  '
$ws.Range("E41").Value = '  +0.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.48'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.29'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.19%  '
$ws.Range("E44").Value = '  -0.82%  '
$ws.Range("D45").Value = '1.775.73'
$ws.Range("E45").Value = '  -0.65%  '
$ws.Range("D46").Value = '0.0₆0114'
$ws.Range("E46").Value = '  +1.73%  '
$ws.Range("E47").Value = '  +0.31%  '
$ws.Range("E48").Value = '  +2.47%  '
$ws.Range("E49").Value = '  -1.84%  '
$ws.Range("E50").Value = '  -0.59%  '
$ws.Range("E51").Value = '  -0.62%  '
